# Updated cryptos list on Mon Jul 22 16:48:40 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and column E (Volume(1h)) hold plain-text strings (e.g. "3.462.93",
# "  -0.90%  "). Force a text number format before writing so Excel does not reinterpret
# the digit-and-dot strings as numbers/dates, then restore the default "Normal" style so
# the cell keeps its original (unstyled) look.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "67.180.68"
Set-TextValue "E2" "  +0.17%  "
Set-TextValue "D3" "3.462.93"
Set-TextValue "E3" "  -0.90%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.09%  "
Set-TextValue "D5" "593.70"
Set-TextValue "E5" "  -0.85%  "
Set-TextValue "D6" "180.32"
Set-TextValue "E6" "  +2.95%  "
Set-TextValue "D7" "0.613"
Set-TextValue "E7" "  +4.29%  "
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.02%  "
Set-TextValue "D9" "3.458.67"
Set-TextValue "E9" "  -0.95%  "
Set-TextValue "D10" "0.140"
Set-TextValue "E10" "  +6.83%  "
Set-TextValue "D11" "6.97"
Set-TextValue "E11" "  -2.43%  "
Set-TextValue "D12" "0.430"
Set-TextValue "E12" "  +0.12%  "
Set-TextValue "D13" "4.053.34"
Set-TextValue "E13" "  -1.12%  "
Set-TextValue "D14" "31.95"
Set-TextValue "E14" "  +2.80%  "
Set-TextValue "E15" "  -0.49%  "
Set-TextValue "D16" "67.035.74"
Set-TextValue "E16" "  -0.06%  "
Set-TextValue "E17" "  -0.71%  "
Set-TextValue "D18" "3.457.82"
Set-TextValue "E18" "  -1.09%  "
Set-TextValue "D19" "6.20"
Set-TextValue "E19" "  -1.24%  "
Set-TextValue "D20" "14.18"
Set-TextValue "E20" "  -2.56%  "
Set-TextValue "D21" "392.85"
Set-TextValue "E21" "  +0.15%  "
Set-TextValue "D22" "7.93"
Set-TextValue "E22" "  -0.82%  "
Set-TextValue "D23" "1.00"
Set-TextValue "E23" "  -0.13%  "
Set-TextValue "D24" "5.77"
Set-TextValue "E24" "  +1.12%  "
Set-TextValue "D25" "0.539"
Set-TextValue "E25" "  +0.42%  "
Set-TextValue "D26" "71.69"
Set-TextValue "E26" "  -2.21%  "
Set-TextValue "D27" "0.0000122"
Set-TextValue "E27" "  +0.36%  "
Set-TextValue "D28" "10.34"
Set-TextValue "E28" "  +1.22%  "
Set-TextValue "D29" "0.175"
Set-TextValue "E29" "  -2.86%  "
Set-TextValue "E30" "  +0.61%  "
Set-TextValue "D31" "6.13"
Set-TextValue "E31" "  +0.61%  "
Set-TextValue "E32" "  -1.33%  "
Set-TextValue "E33" "  -0.59%  "
Set-TextValue "D34" "23.53"
Set-TextValue "E34" "  -0.39%  "
Set-TextValue "D35" "7.32"
Set-TextValue "E35" "  -0.62%  "
Set-TextValue "E36" "  -0.01%  "
Set-TextValue "D37" "1.58"
Set-TextValue "E37" "  -3.24%  "
Set-TextValue "D38" "160.36"
Set-TextValue "E38" "  -1.65%  "
Set-TextValue "D39" "0.876"
Set-TextValue "E39" "  -0.14%  "
Set-TextValue "D40" "2.81"
Set-TextValue "E40" "  +10.81%  "
Set-TextValue "D41" "1.87"
Set-TextValue "E41" "  -3.15%  "
Set-TextValue "D44" "26.16"
Set-TextValue "E44" "  +0.33%  "
Set-TextValue "D45" "0.0720"
Set-TextValue "E45" "  -1.50%  "
Set-TextValue "D46" "2.759.57"
Set-TextValue "E46" "  -1.61%  "
Set-TextValue "D47" "26.21"
Set-TextValue "E47" "  -4.05%  "
Set-TextValue "D48" "41.34"
Set-TextValue "E48" "  -2.78%  "
Set-TextValue "D49" "0.0298"
Set-TextValue "E49" "  -0.81%  "
Set-TextValue "D50" "324.83"
Set-TextValue "E50" "  -3.99%  "
Set-TextValue "E51" "  -2.87%  "

# Filecoin overtook RenderToken in the ranking: rows 42/43 swap coin identity (name + link)
# in addition to refreshed price/volume figures.
Set-TextValue "B42" "Filecoin"
Set-TextValue "C42" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D42" "4.64"
Set-TextValue "E42" "  -0.29%  "
Set-TextValue "B43" "RenderToken"
Set-TextValue "C43" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D43" "6.72"
Set-TextValue "E43" "  -4.42%  "

